$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.543.89'
$ws.Cells.Item(2, 5).Value = '  -7.61%  '
$ws.Cells.Item(3, 4).Value = '2.543.58'
$ws.Cells.Item(3, 5).Value = '  -2.17%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '295.25'
$ws.Cells.Item(5, 5).Value = '  -5.42%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '92.69'
$ws.Cells.Item(6, 5).Value = '  -6.37%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.571'
$ws.Cells.Item(7, 5).Value = '  -4.43%  '
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.544'
$ws.Cells.Item(9, 5).Value = '  -6.09%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '35.24'
$ws.Cells.Item(10, 5).Value = '  -9.36%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0800'
$ws.Cells.Item(11, 5).Value = '  -4.63%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.63'
$ws.Cells.Item(12, 5).Value = '  -6.18%  '
$ws.Cells.Item(13, 4).Value = '2.931.11'
$ws.Cells.Item(13, 5).Value = '  -2.42%  '
$ws.Cells.Item(14, 5).Value = '  -0.05%  '
$ws.Cells.Item(15, 4).Value = '2.533.87'
$ws.Cells.Item(15, 5).Value = '  -2.93%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.858'
$ws.Cells.Item(16, 5).Value = '  -6.19%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '13.99'
$ws.Cells.Item(17, 5).Value = '  -5.58%  '
$ws.Cells.Item(18, 4).Value = '42.551.49'
$ws.Cells.Item(18, 5).Value = '  -7.89%  '
$ws.Cells.Item(19, 2).Value = 'Uniswap'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.59'
$ws.Cells.Item(19, 5).Value = '  -2.02%  '
$ws.Cells.Item(20, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.51'
$ws.Cells.Item(20, 5).Value = '  -2.54%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0961'
$ws.Cells.Item(21, 5).Value = '  -5.50%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '72.12'
$ws.Cells.Item(22, 5).Value = '  +0.39%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '257.23'
$ws.Cells.Item(23, 5).Value = '  -7.33%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.86'
$ws.Cells.Item(24, 5).Value = '  -6.95%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '29.41'
$ws.Cells.Item(25, 5).Value = '  -2.34%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.10'
$ws.Cells.Item(26, 5).Value = '  -5.09%  '
$ws.Cells.Item(27, 5).Value = '  +0.24%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.90'
$ws.Cells.Item(28, 5).Value = '  -7.71%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.11'
$ws.Cells.Item(29, 5).Value = '  -4.44%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '35.73'
$ws.Cells.Item(30, 5).Value = '  -5.57%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '5.86'
$ws.Cells.Item(31, 5).Value = '  -5.94%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '150.09'
$ws.Cells.Item(32, 5).Value = '  -3.74%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.37'
$ws.Cells.Item(33, 5).Value = '  -6.24%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.14'
$ws.Cells.Item(34, 5).Value = '  -3.34%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.72'
$ws.Cells.Item(35, 5).Value = '  -2.91%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0786'
$ws.Cells.Item(36, 5).Value = '  -6.06%  '
$ws.Cells.Item(37, 5).Value = '  -8.66%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '23.94'
$ws.Cells.Item(38, 5).Value = '  +2.00%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.118'
$ws.Cells.Item(39, 5).Value = '  -3.82%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '15.57'
$ws.Cells.Item(40, 5).Value = '  -1.38%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.38'
$ws.Cells.Item(41, 5).Value = '  -5.75%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0306'
$ws.Cells.Item(42, 5).Value = '  -7.30%  '
$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '2.059.40'
$ws.Cells.Item(43, 5).Value = '  -1.82%  '
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.77'
$ws.Cells.Item(44, 5).Value = '  -4.83%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.999'
$ws.Cells.Item(45, 5).Value = '  -0.08%  '
$ws.Cells.Item(46, 5).Value = '  -11.73%  '
$ws.Cells.Item(47, 5).Value = '  +3.34%  '
$ws.Cells.Item(48, 4).Value = '2.787.63'
$ws.Cells.Item(48, 5).Value = '  -2.53%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.68'
$ws.Cells.Item(49, 5).Value = '  -9.25%  '
$ws.Cells.Item(50, 2).Value = 'Stacks'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.68'
$ws.Cells.Item(50, 5).Value = '  -3.69%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '102.77'
$ws.Cells.Item(51, 5).Value = '  -5.47%  '
